$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# The newly inserted column D inherits formatting from column C (Excel default
# "format from left" behavior). Copy the number formatting/style from column E
# (which holds what used to be column D's style) back onto column D so the new
# column matches its row siblings (date style for header rows, number style for
# data rows). Restrict to the contiguous blocks that actually contain data in
# column D/E so we don't materialize the blank separator rows (36, 78) or add a
# stray cell to the section-header rows (37, 79) that never had a column D cell.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting periods figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 3431000
$ws.Range("D9").Value = 1819000
$ws.Range("D10").Value = 1612000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 398000
$ws.Range("D17").Value = 2783000
$ws.Range("D18").Value = 648000
$ws.Range("D20").Value = 26000
$ws.Range("D21").Value = 1072000
$ws.Range("D22").Value = 152000
$ws.Range("D23").Value = 522000
$ws.Range("D24").Value = -1000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 523000
$ws.Range("D27").Value = 485000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -26000
$ws.Range("D33").Value = 485000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 485000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 8000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 309000
$ws.Range("D44").Value = 50000
$ws.Range("D45").Value = 82000
$ws.Range("D46").Value = 449000
$ws.Range("D47").Value = 317000
$ws.Range("D48").Value = 10871000
$ws.Range("D49").Value = 761000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 46000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 12444000
$ws.Range("D57").Value = 292000
$ws.Range("D58").Value = 1149000
$ws.Range("D59").Value = 174000
$ws.Range("D60").Value = 1615000
$ws.Range("D61").Value = 3129000
$ws.Range("D62").Value = 82000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 4864000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 362000
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 7218000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 485000
$ws.Range("D83").Value = 398000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 924000
$ws.Range("D91").Value = -728000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1154000
$ws.Range("D96").Value = -591000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 233000
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 3000
